# Update odds values on Sheet1 to reflect latest Betfair Back/Lay quotes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 3.1
$ws.Range("I2").Value = 2.92

# Row 3
$ws.Range("J3").Value = 3.55
$ws.Range("L3").Value = 1.43
$ws.Range("M3").Value = 1.06
$ws.Range("O3").Value = 1.32
$ws.Range("S3").Value = 3.5

# Row 4
$ws.Range("F4").Value = 2.28
$ws.Range("K4").Value = 4.8

# Row 5
$ws.Range("J5").Value = 3.9

# Row 6
$ws.Range("G6").Value = 3.1
$ws.Range("I6").Value = 2.8
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 1.98
$ws.Range("Q6").Value = 1.82
